$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DB")

# The "ID" primary-key columns (B and C) for each of the four tables
# (rows 5, 23, 38, 52) need their type/constraint text updated so that
# the generated CREATE TABLE column definition reflects an
# auto-incrementing integer primary key instead of a plain DECIMAL key.
$idRows = @(5, 23, 38, 52)

foreach ($r in $idRows) {
    $ws.Range("C$r").Value = "NOT NULL PRIMARY KEY auto_increment "
    $ws.Range("B$r").Value = "int(8)"
}

# Restore the selection state left by the author after making the edit.
$ws.Activate()
$ws.Range("E51:E58").Select()

$wb.Save()
